$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.739.15"
$ws.Range("E2").Value = "  +1.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.83"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.54"
$ws.Range("E5").Value = "  +2.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4710"
$ws.Range("E7").Value = "  +3.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3941"
$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.62"
$ws.Range("E9").Value = "  -1.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08073"
$ws.Range("E10").Value = "  +2.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.029"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.17"
$ws.Range("E12").Value = "  +3.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.886.33"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.981"
$ws.Range("E14").Value = "  +1.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.133"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06717"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001050"
$ws.Range("E18").Value = "  +2.20%  "

$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.16"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.544"
$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.745.68"
$ws.Range("E23").Value = "  +1.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.114.15"
$ws.Range("E26").Value = "  +1.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.53"
$ws.Range("E27").Value = "  +4.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.24"
$ws.Range("E28").Value = "  +1.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.105"
$ws.Range("E29").Value = "  +2.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.591"
$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.02"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9860"
$ws.Range("E32").Value = "  +4.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09478"
$ws.Range("E33").Value = "  +1.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.453"
$ws.Range("E34").Value = "  +0.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.616"
$ws.Range("E35").Value = "  +0.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.367"
$ws.Range("E36").Value = "  +2.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06150"
$ws.Range("E37").Value = "  +1.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02265"
$ws.Range("E38").Value = "  +1.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.230"
$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.120"
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6005"
$ws.Range("E41").Value = "  +1.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1904"
$ws.Range("E42").Value = "  +1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.31"
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.258"
$ws.Range("E44").Value = "  -1.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5732"
$ws.Range("E45").Value = "  +2.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.23"
$ws.Range("E46").Value = "  +1.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.949"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.398"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06909"
$ws.Range("E49").Value = "  +2.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.67"
$ws.Range("E50").Value = "  +6.12%  "

$ws.Range("E51").Value = "  +5.76%  "
